# "Added finance and criterion score in result pages."
#
# The workbook tracks 6 daily result pages ("20 02 2017" .. "27 02 2017").
# Each page scores a stock using a RANDBETWEEN-driven "finance"/"criterion"
# formula (C5/D5, rolled up into G2 via SUM(D2:D5)). Refreshing the pages
# re-rolls those volatile formulas, and the active selection moved around
# as the author clicked through sheets while reviewing the new scores.

$wb = $excel.ActiveWorkbook

# "20 02 2017": cursor left sitting on A5 after reviewing the scores.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A5").Select() | Out-Null

# "23 02 2017": became the active/focused tab, cursor parked at M27.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("M27").Select() | Out-Null

# "27 02 2017": was previously the active tab; cursor/selection moved to
# A20:A24 and it lost focus (the activation below hands focus to ws4).
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A20:A24").Select() | Out-Null

# Make "23 02 2017" the final active sheet/tab (matches activeTab="3").
$ws4.Activate() | Out-Null

# Recompute the workbook so every sheet's finance score (C5), weighted
# criterion score (D5) and total (G2) reflect a fresh roll of the
# RANDBETWEEN-based scoring formulas.
$excel.CalculateFull() | Out-Null
